# Update the "2002 Update GAO" row: the start-of-range date in B2 was
# recorded incorrectly (40575 -> 1-Feb-2011). Correct it to 36923
# (1-Feb-2001), matching the "February 2001" mentioned in D2's label and
# mirroring the corresponding date already used in C10 for the prior block.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 36923

# Leave the active selection on B2 (previously it was left on F10).
$ws.Range("B2").Select()
